$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing staff columns run F:I (工作人員1..4). Extend the roster with two
# more staff columns, 工作人員5 and 工作人員6, in J1 and K1 - matching the
# header style already used by I1.
$ws.Range("J1").Value = "工作人員5"
$ws.Range("K1").Value = "工作人員6"

# Match the column width/formatting already used for the staff columns
# (E:I) by carrying it over to the two new columns (J:K).
$srcWidth = $ws.Range("I1").ColumnWidth
$ws.Range("J1:K1").ColumnWidth = $srcWidth

# The duplicate-value conditional formatting that covers the staff columns
# (previously E1:I1048576) needs to widen to include the new columns.
$cf = $ws.Range("E1:I1048576").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("E1:K1048576"))

# Leave the selection on A2, as in the saved workbook.
[void]$ws.Range("A2").Select()
